$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.645.44"
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("D3").Value = "2.696.50"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'524.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").Value = "'147.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.578"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").Value = "2.716.81"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").Value = "'6.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.69%  "
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("D14").Value = "3.160.70"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "60.779.14"
$ws.Range("E15").Value = "  +3.26%  "
$ws.Range("D16").Value = "'21.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.772.80"
$ws.Range("E17").Value = "  +3.79%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000139"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").Value = "'351.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "'4.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "'10.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("D22").Value = "'6.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.00%  "
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "'63.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.82%  "
$ws.Range("D25").Value = "'0.425"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("E26").Value = "  +5.33%  "
$ws.Range("D27").Value = "'0.994"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "0.0₃0820"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").Value = "'7.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("D30").Value = "'6.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.64%  "
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  +1.95%  "
$ws.Range("D33").Value = "'19.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").Value = "'147.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  +7.43%  "
$ws.Range("E36").Value = "  +9.86%  "
$ws.Range("D37").Value = "'0.954"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.63%  "
$ws.Range("E38").Value = "  +11.07%  "
$ws.Range("D39").Value = "'0.880"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.03%  "
$ws.Range("D40").Value = "'36.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'3.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").Value = "'283.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("D43").Value = "'0.615"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'20.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.0990"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "2.142.06"
$ws.Range("E47").Value = "  +7.50%  "
$ws.Range("D48").Value = "'0.0539"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("D49").Value = "'4.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.89%  "
$ws.Range("E50").Value = "  +2.38%  "
$ws.Range("D51").Value = "'10.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.02%  "
